# Weekly update: two new price rows for "Fruta, Macroferia Regional de Talca - Pera".
# Insert two rows at the top of the data block (row 277) so the rest of the
# existing rows shift down by two, then populate the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(277).Resize(2).Insert()

# New row 277
$ws.Cells.Item(277, 1).Value = 5
$ws.Cells.Item(277, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(277, 3).Value = "Maule"
$ws.Cells.Item(277, 4).Value = 44461
$ws.Cells.Item(277, 5).Value = 7
$ws.Cells.Item(277, 6).Value = "Fruta"
$ws.Cells.Item(277, 7).Value = 100104
$ws.Cells.Item(277, 8).Value = "Frutos de pepita"
$ws.Cells.Item(277, 9).Value = 100104005
$ws.Cells.Item(277, 10).Value = "Pera"
$ws.Cells.Item(277, 11).Value = "Packham's Triumph"
$ws.Cells.Item(277, 12).Value = "Especial"
$ws.Cells.Item(277, 13).Value = 360
$ws.Cells.Item(277, 14).Value = 9000
$ws.Cells.Item(277, 15).Value = 9000
$ws.Cells.Item(277, 16).Value = 9000
$ws.Cells.Item(277, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(277, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(277, 19).Value = 500
$ws.Cells.Item(277, 20).Value = 18

# New row 278
$ws.Cells.Item(278, 1).Value = 5
$ws.Cells.Item(278, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(278, 3).Value = "Maule"
$ws.Cells.Item(278, 4).Value = 44461
$ws.Cells.Item(278, 5).Value = 7
$ws.Cells.Item(278, 6).Value = "Fruta"
$ws.Cells.Item(278, 7).Value = 100104
$ws.Cells.Item(278, 8).Value = "Frutos de pepita"
$ws.Cells.Item(278, 9).Value = 100104005
$ws.Cells.Item(278, 10).Value = "Pera"
$ws.Cells.Item(278, 11).Value = "Packham's Triumph"
$ws.Cells.Item(278, 12).Value = "Primera"
$ws.Cells.Item(278, 13).Value = 200
$ws.Cells.Item(278, 14).Value = 8000
$ws.Cells.Item(278, 15).Value = 8000
$ws.Cells.Item(278, 16).Value = 8000
$ws.Cells.Item(278, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(278, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(278, 19).Value = 444
$ws.Cells.Item(278, 20).Value = 18
